$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Simple value updates (Price / Volume(1h)) ---

Set-TextCell "D2" "37.118.98"
$ws.Range("E2").Value = "  -0.85%  "

Set-TextCell "D3" "2.014.70"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  +0.21%  "

Set-TextCell "D5" "226.80"
$ws.Range("E5").Value = "  -0.73%  "

Set-TextCell "D6" "0.608"
$ws.Range("E6").Value = "  -0.62%  "

Set-TextCell "D8" "55.24"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("E9").Value = "  -3.31%  "

Set-TextCell "D10" "0.0778"
$ws.Range("E10").Value = "  -3.35%  "

$ws.Range("E11").Value = "  -4.03%  "

Set-TextCell "D12" "2.311.77"
$ws.Range("E12").Value = "  -1.74%  "

Set-TextCell "D13" "14.06"
$ws.Range("E13").Value = "  -3.02%  "

Set-TextCell "D14" "19.79"
$ws.Range("E14").Value = "  -3.74%  "

Set-TextCell "D15" "0.737"
$ws.Range("E15").Value = "  -2.23%  "

$ws.Range("E16").Value = "  -1.70%  "

Set-TextCell "D17" "2.002.70"
$ws.Range("E17").Value = "  -2.12%  "

Set-TextCell "D18" "36.984.36"
$ws.Range("E18").Value = "  -0.93%  "

Set-TextCell "D19" "6.22"
$ws.Range("E19").Value = "  +2.61%  "

Set-TextCell "D20" "68.89"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("E21").Value = "  -4.22%  "

Set-TextCell "D22" "222.81"
$ws.Range("E22").Value = "  -1.21%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  +2.71%  "

Set-TextCell "D25" "2.18"
$ws.Range("E25").Value = "  -4.51%  "

Set-TextCell "D26" "164.74"
$ws.Range("E26").Value = "  -2.11%  "

Set-TextCell "D27" "8.97"
$ws.Range("E27").Value = "  -5.34%  "

$ws.Range("E28").Value = "  -2.05%  "

Set-TextCell "D29" "18.61"
$ws.Range("E29").Value = "  -1.47%  "

$ws.Range("E30").Value = "  -5.06%  "

$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("E32").Value = "  -2.69%  "

$ws.Range("E33").Value = "  -1.73%  "

$ws.Range("E34").Value = "  -1.65%  "

# --- Row 35/36 swap (LidoDAOToken <-> WEMIXToken) with refreshed data ---

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D35" "1.90"
$ws.Range("E35").Value = "  +4.31%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D36" "2.34"
$ws.Range("E36").Value = "  -2.14%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -1.48%  "

Set-TextCell "D39" "5.39"
$ws.Range("E39").Value = "  -0.98%  "

Set-TextCell "D40" "1.463.46"
$ws.Range("E40").Value = "  -2.24%  "

Set-TextCell "D41" "0.0212"
$ws.Range("E41").Value = "  -3.95%  "

Set-TextCell "D42" "94.37"
$ws.Range("E42").Value = "  -1.57%  "

Set-TextCell "D43" "0.0908"
$ws.Range("E43").Value = "  -2.88%  "

$ws.Range("E44").Value = "  -3.80%  "

# --- Row 45/46 swap (FTXToken <-> InjectiveProtocol) with refreshed data ---

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D45" "15.92"
$ws.Range("E45").Value = "  -4.96%  "

$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D46" "4.12"
$ws.Range("E46").Value = "  +9.44%  "

Set-TextCell "D47" "1.11"
$ws.Range("E47").Value = "  -2.48%  "

$ws.Range("E48").Value = "  -1.61%  "

$ws.Range("E49").Value = "  -2.43%  "

Set-TextCell "D50" "2.91"
$ws.Range("E50").Value = "  -0.51%  "

Set-TextCell "D51" "2.198.90"
$ws.Range("E51").Value = "  -1.79%  "
